{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer paragraphs,\n// along with the blank paragraph that precedes them, while leaving the\n// \"LOQ4064: ...\" requirement paragraph and the remaining blank / page-break\n// paragraphs that follow untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the footer paragraphs by their text content so the script is\n// resilient to the exact paragraph count/position.\nlet jupiterIdx = -1;\nlet copyrightIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (jupiterIdx === -1 && t.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIdx = i;\n  }\n  if (copyrightIdx === -1 && t.indexOf(\"Contact: luizeleno@usp.br\") !== -1) {\n    copyrightIdx = i;\n  }\n}\n\nif (jupiterIdx !== -1 && copyrightIdx !== -1) {\n  // The blank paragraph immediately before the \"Ver no Jupiter\" paragraph\n  // is also removed (it only separated the footer block from the rest).\n  let blankIdx = -1;\n  if (jupiterIdx - 1 >= 0 && items[jupiterIdx - 1].text === \"\") {\n    blankIdx = jupiterIdx - 1;\n  }\n\n  // Delete from the end backwards so earlier indices stay valid.\n  items[copyrightIdx].delete();\n  items[jupiterIdx].delete();\n  if (blankIdx !== -1) {\n    items[blankIdx].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"(c) 2020 ...\" footer paragraphs,\n# along with the blank paragraph that precedes them, while leaving the\n# \"LOQ4064: ...\" requirement paragraph and the remaining blank / page-break\n# paragraphs that follow untouched.\n$d = $word.ActiveDocument\n\n$jupiterIdx = -1\n$copyrightIdx = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($jupiterIdx -eq -1 -and $t -like \"*Ver no Jupiter*\") {\n        $jupiterIdx = $i\n    }\n    if ($copyrightIdx -eq -1 -and $t -like \"*Contact: luizeleno@usp.br*\") {\n        $copyrightIdx = $i\n    }\n}\n\nif ($jupiterIdx -ne -1 -and $copyrightIdx -ne -1) {\n    # The blank paragraph immediately before the \"Ver no Jupiter\" paragraph\n    # only separated the footer block from the rest, so remove it too.\n    $blankIdx = -1\n    if ($jupiterIdx - 1 -ge 1) {\n        $prevText = $d.Paragraphs.Item($jupiterIdx - 1).Range.Text\n        if ($prevText -eq \"`r\") {\n            $blankIdx = $jupiterIdx - 1\n        }\n    }\n\n    # Delete from the end backwards so earlier indices stay valid.\n    $d.Paragraphs.Item($copyrightIdx).Range.Delete() | Out-Null\n    $d.Paragraphs.Item($jupiterIdx).Range.Delete() | Out-Null\n    if ($blankIdx -ne -1) {\n        $d.Paragraphs.Item($blankIdx).Range.Delete() | Out-Null\n    }\n}\n"}
